# FAQ: remove incomplete Q/A
#
# The document had an incomplete "Legal" section appended at the very
# end of the FAQ: a "Legal" Heading2, a single Heading3 question ("As
# an operator of a gateway, am I liable for the use of the gateway by
# "bad actors"?"), and a FirstParagraph answer that was never actually
# written (just the placeholder text "--"). Remove all three
# paragraphs, so the FAQ ends with the previous answer about gateways
# and user privacy.

$d = $word.ActiveDocument

$legalIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Range.ParagraphFormat.Style.NameLocal
    if ($styleName -eq "Heading 2" -and $p.Range.Text -match "^Legal\r?$") {
        $legalIndex = $i
        break
    }
}

if ($legalIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($legalIndex)
    $endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
